# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    # Force the value to be stored as text even if it looks numeric,
    # then restore the default "Normal" style so no formatting residue is left behind.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$ws.Range("D2").Value = "56.857.18"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "2.969.18"
$ws.Range("E3").Value = "  -1.06%  "

Set-TextValue $ws.Range("D4") "0.997"
$ws.Range("E4").Value = "  -0.26%  "

Set-TextValue $ws.Range("D5") "496.83"
$ws.Range("E5").Value = "  -3.05%  "

Set-TextValue $ws.Range("D6") "136.58"
$ws.Range("E6").Value = "  -1.48%  "

Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.04%  "

Set-TextValue $ws.Range("D8") "0.425"
$ws.Range("E8").Value = "  -2.23%  "

Set-TextValue $ws.Range("D9") "7.29"
$ws.Range("E9").Value = "  -3.02%  "

Set-TextValue $ws.Range("D10") "0.106"
$ws.Range("E10").Value = "  -1.60%  "

Set-TextValue $ws.Range("D11") "0.355"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").Value = "3.475.96"
$ws.Range("E12").Value = "  -1.12%  "

$ws.Range("E13").Value = "  -1.60%  "

Set-TextValue $ws.Range("D14") "25.66"
$ws.Range("E14").Value = "  +0.02%  "

Set-TextValue $ws.Range("D15") "0.0000156"
$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("D16").Value = "56.852.15"
$ws.Range("E16").Value = "  +0.03%  "

Set-TextValue $ws.Range("D17") "6.04"
$ws.Range("E17").Value = "  +2.18%  "

$ws.Range("D18").Value = "2.965.94"
$ws.Range("E18").Value = "  -1.20%  "

Set-TextValue $ws.Range("D19") "12.55"
$ws.Range("E19").Value = "  +0.17%  "

Set-TextValue $ws.Range("D20") "7.78"
$ws.Range("E20").Value = "  -0.75%  "

Set-TextValue $ws.Range("D21") "318.44"
$ws.Range("E21").Value = "  -2.54%  "

Set-TextValue $ws.Range("D22") "0.998"
$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("E23").Value = "  -0.83%  "

Set-TextValue $ws.Range("D24") "0.484"
$ws.Range("E24").Value = "  -0.17%  "

Set-TextValue $ws.Range("D25") "62.95"
$ws.Range("E25").Value = "  -0.50%  "

$ws.Range("E26").Value = "  +0.37%  "

Set-TextValue $ws.Range("D27") "0.162"
$ws.Range("E27").Value = "  -5.22%  "

$ws.Range("D28").Value = "0.0₃0885"
$ws.Range("E28").Value = "  -3.22%  "

Set-TextValue $ws.Range("D29") "6.50"
$ws.Range("E29").Value = "  -1.75%  "

Set-TextValue $ws.Range("D30") "7.05"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("E31").Value = "  -2.80%  "

Set-TextValue $ws.Range("D32") "20.10"
$ws.Range("E32").Value = "  -2.20%  "

Set-TextValue $ws.Range("D33") "1.15"
$ws.Range("E33").Value = "  -6.80%  "

Set-TextValue $ws.Range("D34") "154.49"
$ws.Range("E34").Value = "  -1.73%  "

Set-TextValue $ws.Range("D35") "4.58"
$ws.Range("E35").Value = "  +0.27%  "

Set-TextValue $ws.Range("D36") "5.70"
$ws.Range("E36").Value = "  +0.10%  "

Set-TextValue $ws.Range("D37") "1.24"
$ws.Range("E37").Value = "  -2.29%  "

Set-TextValue $ws.Range("D38") "23.89"
$ws.Range("E38").Value = "  +0.03%  "

Set-TextValue $ws.Range("D39") "0.0662"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").Value = "2.995.70"
$ws.Range("E40").Value = "  -1.25%  "

Set-TextValue $ws.Range("D41") "37.43"
$ws.Range("E41").Value = "  +0.79%  "

Set-TextValue $ws.Range("D42") "0.998"
$ws.Range("E42").Value = "  -0.15%  "

Set-TextValue $ws.Range("D43") "3.70"
$ws.Range("E43").Value = "  +0.90%  "

Set-TextValue $ws.Range("D44") "0.637"
$ws.Range("E44").Value = "  -1.86%  "

$ws.Range("D45").Value = "2.185.31"
$ws.Range("E45").Value = "  -4.45%  "

Set-TextValue $ws.Range("D46") "1.37"
$ws.Range("E46").Value = "  -3.05%  "

Set-TextValue $ws.Range("D47") "0.942"
$ws.Range("E47").Value = "  -5.77%  "

Set-TextValue $ws.Range("D48") "5.91"
$ws.Range("E48").Value = "  +0.65%  "

Set-TextValue $ws.Range("D49") "0.0234"
$ws.Range("E49").Value = "  -2.29%  "

Set-TextValue $ws.Range("D50") "19.03"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D51") "1.77"
$ws.Range("E51").Value = "  -9.31%  "
